# Update TPM-derived NATMI Hbegf-Cd44 LR-pair metrics on Sheet1.
# Ligand/receptor expression summaries and derived edge weights/specificities
# were recomputed with new TPM values; counts (columns E,F,K,L) are unchanged.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 8.280371333333333
$ws.Range("H2").Value = 24.841114
$ws.Range("I2").Value = 0.2946400644635011
$ws.Range("J2").Value = 0.3116548779253407
$ws.Range("M2").Value = 16.14072933333334
$ws.Range("N2").Value = 48.42218800000001
$ws.Range("O2").Value = 0.03423048004954622
$ws.Range("P2").Value = 0.03634868370049611
$ws.Range("Q2").Value = 133.6512324708258
$ws.Range("R2").Value = 1202.861092237432
$ws.Range("S2").Value = 0.01008567084841489
$ws.Range("T2").Value = 0.01132824458142494
$ws.Range("G3").Value = 8.280371333333333
$ws.Range("H3").Value = 24.841114
$ws.Range("I3").Value = 0.2946400644635011
$ws.Range("J3").Value = 0.3116548779253407
$ws.Range("O3").Value = 0.1719151703242873
$ws.Range("P3").Value = 0.1825533892714798
$ws.Range("Q3").Value = 671.234360751465
$ws.Range("R3").Value = 6041.109246763187
$ws.Range("S3").Value = 0.05065309686660179
$ws.Range("T3").Value = 0.05689365424826026
$ws.Range("G4").Value = 8.280371333333333
$ws.Range("H4").Value = 24.841114
$ws.Range("I4").Value = 0.2946400644635011
$ws.Range("J4").Value = 0.3116548779253407
$ws.Range("M4").Value = 168.70371
$ws.Range("N4").Value = 506.11113
$ws.Range("O4").Value = 0.3577786889414888
$ws.Range("P4").Value = 0.3799182594076638
$ws.Range("Q4").Value = 1396.92936411098
$ws.Range("R4").Value = 12572.36427699882
$ws.Range("S4").Value = 0.1054159359733872
$ws.Range("T4").Value = 0.1184033787573034
$ws.Range("G5").Value = 8.280371333333333
$ws.Range("H5").Value = 24.841114
$ws.Range("I5").Value = 0.2946400644635011
$ws.Range("J5").Value = 0.3116548779253407
$ws.Range("M5").Value = 82.43477250000001
$ws.Range("N5").Value = 164.869545
$ws.Range("O5").Value = 0.1748236883957081
$ws.Range("P5").Value = 0.1237612588479007
$ws.Range("Q5").Value = 682.5905270788551
$ws.Range("R5").Value = 4095.54316247313
$ws.Range("S5").Value = 0.05151006281865845
$ws.Range("T5").Value = 0.03857080001812899
$ws.Range("G6").Value = 8.280371333333333
$ws.Range("H6").Value = 24.841114
$ws.Range("I6").Value = 0.2946400644635011
$ws.Range("J6").Value = 0.3116548779253407
$ws.Range("M6").Value = 123.1883796666667
$ws.Range("N6").Value = 369.565139
$ws.Range("O6").Value = 0.2612519722889696
$ws.Range("P6").Value = 0.2774184087724594
$ws.Range("Q6").Value = 1020.045527591649
$ws.Range("R6").Value = 9180.409748324846
$ws.Range("S6").Value = 0.0769752979564388
$ws.Range("T6").Value = 0.08645880032022313
$ws.Range("I7").Value = 0.405746032520008
$ws.Range("J7").Value = 0.4291769704298953
$ws.Range("M7").Value = 16.14072933333334
$ws.Range("N7").Value = 48.42218800000001
$ws.Range("O7").Value = 0.03423048004954622
$ws.Range("P7").Value = 0.03634868370049611
$ws.Range("Q7").Value = 184.0498420172062
$ws.Range("R7").Value = 1656.448578154856
$ws.Range("S7").Value = 0.01388888147135867
$ws.Range("T7").Value = 0.01560001794969344
$ws.Range("I8").Value = 0.405746032520008
$ws.Range("J8").Value = 0.4291769704298953
$ws.Range("O8").Value = 0.1719151703242873
$ws.Range("P8").Value = 0.1825533892714798
$ws.Range("S8").Value = 0.069753898289081
$ws.Range("T8").Value = 0.07834771054924307
$ws.Range("I9").Value = 0.405746032520008
$ws.Range("J9").Value = 0.4291769704298953
$ws.Range("M9").Value = 168.70371
$ws.Range("N9").Value = 506.11113
$ws.Range("O9").Value = 0.3577786889414888
$ws.Range("P9").Value = 0.3799182594076638
$ws.Range("Q9").Value = 1923.69815093134
$ws.Range("R9").Value = 17313.28335838206
$ws.Range("S9").Value = 0.1451672835582191
$ws.Range("T9").Value = 0.1630521675835802
$ws.Range("I10").Value = 0.405746032520008
$ws.Range("J10").Value = 0.4291769704298953
$ws.Range("M10").Value = 82.43477250000001
$ws.Range("N10").Value = 164.869545
$ws.Range("O10").Value = 0.1748236883957081
$ws.Range("P10").Value = 0.1237612588479007
$ws.Range("Q10").Value = 939.988927514965
$ws.Range("R10").Value = 5639.93356508979
$ws.Range("S10").Value = 0.07093401795707271
$ws.Range("T10").Value = 0.05311548212893209
$ws.Range("I11").Value = 0.405746032520008
$ws.Range("J11").Value = 0.4291769704298953
$ws.Range("M11").Value = 123.1883796666667
$ws.Range("N11").Value = 369.565139
$ws.Range("O11").Value = 0.2612519722889696
$ws.Range("P11").Value = 0.2774184087724594
$ws.Range("Q11").Value = 1404.695001556246
$ws.Range("R11").Value = 12642.25501400622
$ws.Range("S11").Value = 0.1060019512442765
$ws.Range("T11").Value = 0.1190615922184464
$ws.Range("G12").Value = 1.864050333333333
$ws.Range("H12").Value = 5.592150999999999
$ws.Range("I12").Value = 0.06632841551025578
$ws.Range("J12").Value = 0.07015873512134246
$ws.Range("M12").Value = 16.14072933333334
$ws.Range("N12").Value = 48.42218800000001
$ws.Range("O12").Value = 0.03423048004954622
$ws.Range("P12").Value = 0.03634868370049611
$ws.Range("Q12").Value = 30.08713189404311
$ws.Range("R12").Value = 270.784187046388
$ws.Range("S12").Value = 0.002270453503841823
$ws.Range("T12").Value = 0.002550177671752565
$ws.Range("G13").Value = 1.864050333333333
$ws.Range("H13").Value = 5.592150999999999
$ws.Range("I13").Value = 0.06632841551025578
$ws.Range("J13").Value = 0.07015873512134246
$ws.Range("O13").Value = 0.1719151703242873
$ws.Range("P13").Value = 0.1825533892714798
$ws.Range("Q13").Value = 151.1061018322554
$ws.Range("R13").Value = 1359.954916490299
$ws.Range("S13").Value = 0.01140286084978572
$ws.Range("T13").Value = 0.01280771488340108
$ws.Range("G14").Value = 1.864050333333333
$ws.Range("H14").Value = 5.592150999999999
$ws.Range("I14").Value = 0.06632841551025578
$ws.Range("J14").Value = 0.07015873512134246
$ws.Range("M14").Value = 168.70371
$ws.Range("N14").Value = 506.11113
$ws.Range("O14").Value = 0.3577786889414888
$ws.Range("P14").Value = 0.3799182594076638
$ws.Range("Q14").Value = 314.4722068600699
$ws.Range("R14").Value = 2830.24986174063
$ws.Range("S14").Value = 0.02373089354082562
$ws.Range("T14").Value = 0.02665458452954376
$ws.Range("G15").Value = 1.864050333333333
$ws.Range("H15").Value = 5.592150999999999
$ws.Range("I15").Value = 0.06632841551025578
$ws.Range("J15").Value = 0.07015873512134246
$ws.Range("M15").Value = 82.43477250000001
$ws.Range("N15").Value = 164.869545
$ws.Range("O15").Value = 0.1748236883957081
$ws.Range("P15").Value = 0.1237612588479007
$ws.Range("Q15").Value = 153.6625651568825
$ws.Range("R15").Value = 921.975390941295
$ws.Range("S15").Value = 0.01159577824494601
$ws.Range("T15").Value = 0.008682933377793765
$ws.Range("G16").Value = 1.864050333333333
$ws.Range("H16").Value = 5.592150999999999
$ws.Range("I16").Value = 0.06632841551025578
$ws.Range("J16").Value = 0.07015873512134246
$ws.Range("M16").Value = 123.1883796666667
$ws.Range("N16").Value = 369.565139
$ws.Range("O16").Value = 0.2612519722889696
$ws.Range("P16").Value = 0.2774184087724594
$ws.Range("Q16").Value = 229.6293401804432
$ws.Range("R16").Value = 2066.664061623989
$ws.Range("S16").Value = 0.0173284293708566
$ws.Range("T16").Value = 0.01946332465885129
$ws.Range("G17").Value = 4.6029105
$ws.Range("H17").Value = 9.205821
$ws.Range("I17").Value = 0.1637851482553954
$ws.Range("J17").Value = 0.1154955860658076
$ws.Range("M17").Value = 16.14072933333334
$ws.Range("N17").Value = 48.42218800000001
$ws.Range("O17").Value = 0.03423048004954622
$ws.Range("P17").Value = 0.03634868370049611
$ws.Range("Q17").Value = 74.29433252605801
$ws.Range("R17").Value = 445.7659951563481
$ws.Range("S17").Value = 0.005606444249768283
$ws.Range("T17").Value = 0.004198112526709467
$ws.Range("G18").Value = 4.6029105
$ws.Range("H18").Value = 9.205821
$ws.Range("I18").Value = 0.1637851482553954
$ws.Range("J18").Value = 0.1154955860658076
$ws.Range("O18").Value = 0.1719151703242873
$ws.Range("P18").Value = 0.1825533892714798
$ws.Range("Q18").Value = 373.1271899155215
$ws.Range("R18").Value = 2238.763139493129
$ws.Range("S18").Value = 0.02815715165891496
$ws.Range("T18").Value = 0.02108411068220908
$ws.Range("G19").Value = 4.6029105
$ws.Range("H19").Value = 9.205821
$ws.Range("I19").Value = 0.1637851482553954
$ws.Range("J19").Value = 0.1154955860658076
$ws.Range("M19").Value = 168.70371
$ws.Range("N19").Value = 506.11113
$ws.Range("O19").Value = 0.3577786889414888
$ws.Range("P19").Value = 0.3799182594076638
$ws.Range("Q19").Value = 776.5280781479551
$ws.Range("R19").Value = 4659.16846888773
$ws.Range("S19").Value = 0.05859883561090274
$ws.Range("T19").Value = 0.04387888202738966
$ws.Range("G20").Value = 4.6029105
$ws.Range("H20").Value = 9.205821
$ws.Range("I20").Value = 0.1637851482553954
$ws.Range("J20").Value = 0.1154955860658076
$ws.Range("M20").Value = 82.43477250000001
$ws.Range("N20").Value = 164.869545
$ws.Range("O20").Value = 0.1748236883957081
$ws.Range("P20").Value = 0.1237612588479007
$ws.Range("Q20").Value = 379.4398799053613
$ws.Range("R20").Value = 1517.759519621445
$ws.Range("S20").Value = 0.0286335237224461
$ws.Range("T20").Value = 0.01429387912288041
$ws.Range("G21").Value = 4.6029105
$ws.Range("H21").Value = 9.205821
$ws.Range("I21").Value = 0.1637851482553954
$ws.Range("J21").Value = 0.1154955860658076
$ws.Range("M21").Value = 123.1883796666667
$ws.Range("N21").Value = 369.565139
$ws.Range("O21").Value = 0.2612519722889696
$ws.Range("P21").Value = 0.2774184087724594
$ws.Range("Q21").Value = 567.0250862456865
$ws.Range("R21").Value = 3402.150517474119
$ws.Range("S21").Value = 0.04278919301336333
$ws.Range("T21").Value = 0.03204060170661899
$ws.Range("G22").Value = 1.953192
$ws.Range("H22").Value = 5.859576
$ws.Range("I22").Value = 0.0695003392508397
$ws.Range("J22").Value = 0.0735138304576138
$ws.Range("M22").Value = 16.14072933333334
$ws.Range("N22").Value = 48.42218800000001
$ws.Range("O22").Value = 0.03423048004954622
$ws.Range("P22").Value = 0.03634868370049611
$ws.Range("Q22").Value = 31.525943408032
$ws.Range("R22").Value = 283.733490672288
$ws.Range("S22").Value = 0.002379029976162563
$ws.Range("T22").Value = 0.002672130970915701
$ws.Range("G23").Value = 1.953192
$ws.Range("H23").Value = 5.859576
$ws.Range("I23").Value = 0.0695003392508397
$ws.Range("J23").Value = 0.0735138304576138
$ws.Range("O23").Value = 0.1719151703242873
$ws.Range("P23").Value = 0.1825533892714798
$ws.Range("Q23").Value = 158.332220955736
$ws.Range("R23").Value = 1424.989988601624
$ws.Range("S23").Value = 0.01194816265990386
$ws.Range("T23").Value = 0.01342019890836634
$ws.Range("G24").Value = 1.953192
$ws.Range("H24").Value = 5.859576
$ws.Range("I24").Value = 0.0695003392508397
$ws.Range("J24").Value = 0.0735138304576138
$ws.Range("M24").Value = 168.70371
$ws.Range("N24").Value = 506.11113
$ws.Range("O24").Value = 0.3577786889414888
$ws.Range("P24").Value = 0.3799182594076638
$ws.Range("Q24").Value = 329.51073674232
$ws.Range("R24").Value = 2965.59663068088
$ws.Range("S24").Value = 0.02486574025815412
$ws.Range("T24").Value = 0.02792924650984674
$ws.Range("G25").Value = 1.953192
$ws.Range("H25").Value = 5.859576
$ws.Range("I25").Value = 0.0695003392508397
$ws.Range("J25").Value = 0.0735138304576138
$ws.Range("M25").Value = 82.43477250000001
$ws.Range("N25").Value = 164.869545
$ws.Range("O25").Value = 0.1748236883957081
$ws.Range("P25").Value = 0.1237612588479007
$ws.Range("Q25").Value = 161.01093816882
$ws.Range("R25").Value = 966.06562901292
$ws.Range("S25").Value = 0.0121503056525848
$ws.Range("T25").Value = 0.009098164200165426
$ws.Range("G26").Value = 1.953192
$ws.Range("H26").Value = 5.859576
$ws.Range("I26").Value = 0.0695003392508397
$ws.Range("J26").Value = 0.0735138304576138
$ws.Range("M26").Value = 123.1883796666667
$ws.Range("N26").Value = 369.565139
$ws.Range("O26").Value = 0.2612519722889696
$ws.Range("P26").Value = 0.2774184087724594
$ws.Range("Q26").Value = 240.610557657896
$ws.Range("R26").Value = 2165.495018921064
$ws.Range("S26").Value = 0.01815710070403436
$ws.Range("T26").Value = 0.02039408986831958
